$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the empty, unused second sheet.
$null = $wb.Worksheets.Item("Sheet2").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new "site" column right before the "species" column (old column E),
# shifting species/lat/lng/x/y/z one column to the right.
$null = $ws.Columns.Item(5).Insert()

$ws.Range("E1").Value = "site"
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 5).Value = "Mojave"
}

# Restore the selection to where the author last left it.
$null = $ws.Range("C39").Select()
